# Review mark-ups: flip the sign of every Z-score value in column M
# ("category2") for data rows 2 through 204 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 204
$col      = 13  # Column M

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = -1 * $current
    }
}
